# Auto-generated edit script: apply COVID data update for 'paises.xlsx'
# Updates country stats (new Paraguay/Chad/Liberia/Nueva Caledonia/Montserrat rows
# inserted in sorted position, pushing subsequent rows down) and refreshes the
# 'Datos actualizados' timestamp string, matching the upstream diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 8 de Mayo de 2020 a las 21:34'

# Row 4
$ws.Cells.Item(4, 2).Value = 1309657
$ws.Cells.Item(4, 3).Value = 17034
$ws.Cells.Item(4, 5).Value = 1010710
$ws.Cells.Item(4, 7).Value = 1036
$ws.Cells.Item(4, 8).Value = 77964

# Row 10
$ws.Cells.Item(10, 2).Value = 170114
$ws.Cells.Item(10, 3).Value = 684
$ws.Cells.Item(10, 5).Value = 20969
$ws.Cells.Item(10, 7).Value = 53
$ws.Cells.Item(10, 8).Value = 7445

# Row 48
$ws.Cells.Item(48, 6).Value = 61

# Row 106
$ws.Cells.Item(106, 2).Value = 773
$ws.Cells.Item(106, 3).Value = 8
$ws.Cells.Item(106, 4).Value = 461
$ws.Cells.Item(106, 5).Value = 306

# Row 118
$ws.Cells.Item(118, 1).Value = 'Paraguay'
$ws.Cells.Item(118, 2).Value = 563
$ws.Cells.Item(118, 3).Value = 101
$ws.Cells.Item(118, 4).Value = 152
$ws.Cells.Item(118, 5).Value = 401
$ws.Cells.Item(118, 6).Value = 9
$ws.Cells.Item(118, 8).Value = 10

# Row 119
$ws.Cells.Item(119, 1).Value = 'Tayikistan'
$ws.Cells.Item(119, 2).Value = 522
$ws.Cells.Item(119, 3).Value = 61
$ws.Cells.Item(119, 4).Value = 0
$ws.Cells.Item(119, 5).Value = 510
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 8).Value = 12

# Row 120
$ws.Cells.Item(120, 1).Value = 'Jordania'
$ws.Cells.Item(120, 2).Value = 508
$ws.Cells.Item(120, 3).Value = 14
$ws.Cells.Item(120, 4).Value = 385
$ws.Cells.Item(120, 5).Value = 114
$ws.Cells.Item(120, 6).Value = 5
$ws.Cells.Item(120, 8).Value = 9

# Row 121
$ws.Cells.Item(121, 1).Value = 'Gabon'
$ws.Cells.Item(121, 2).Value = 504
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 110
$ws.Cells.Item(121, 5).Value = 386
$ws.Cells.Item(121, 6).Value = 1
$ws.Cells.Item(121, 8).Value = 8

# Row 122
$ws.Cells.Item(122, 1).Value = 'Malta'
$ws.Cells.Item(122, 2).Value = 489
$ws.Cells.Item(122, 3).Value = 3
$ws.Cells.Item(122, 4).Value = 419
$ws.Cells.Item(122, 5).Value = 65
$ws.Cells.Item(122, 8).Value = 5

# Row 123
$ws.Cells.Item(123, 1).Value = 'Jamaica'
$ws.Cells.Item(123, 2).Value = 488
$ws.Cells.Item(123, 3).Value = 10
$ws.Cells.Item(123, 4).Value = 58
$ws.Cells.Item(123, 5).Value = 421
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 8).Value = 9

# Row 124
$ws.Cells.Item(124, 1).Value = 'Tanzania'
$ws.Cells.Item(124, 2).Value = 480
$ws.Cells.Item(124, 4).Value = 167
$ws.Cells.Item(124, 5).Value = 297
$ws.Cells.Item(124, 6).Value = 7
$ws.Cells.Item(124, 8).Value = 16

# Row 135
$ws.Cells.Item(135, 2).Value = 273
$ws.Cells.Item(135, 3).Value = 2
$ws.Cells.Item(135, 4).Value = 136
$ws.Cells.Item(135, 5).Value = 137

# Row 136
$ws.Cells.Item(136, 1).Value = 'Republica del Chad'
$ws.Cells.Item(136, 2).Value = 260
$ws.Cells.Item(136, 3).Value = 7
$ws.Cells.Item(136, 4).Value = 50
$ws.Cells.Item(136, 5).Value = 182
$ws.Cells.Item(136, 8).Value = 28

# Row 137
$ws.Cells.Item(137, 1).Value = 'Sierra Leona'
$ws.Cells.Item(137, 2).Value = 257
$ws.Cells.Item(137, 3).Value = 26
$ws.Cells.Item(137, 4).Value = 54
$ws.Cells.Item(137, 5).Value = 186
$ws.Cells.Item(137, 7).Value = 1
$ws.Cells.Item(137, 8).Value = 17

# Row 141
$ws.Cells.Item(141, 1).Value = 'Liberia'
$ws.Cells.Item(141, 2).Value = 199
$ws.Cells.Item(141, 3).Value = 10
$ws.Cells.Item(141, 4).Value = 79
$ws.Cells.Item(141, 5).Value = 100
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 8).Value = 20

# Row 142
$ws.Cells.Item(142, 1).Value = 'Etiopia'
$ws.Cells.Item(142, 2).Value = 194
$ws.Cells.Item(142, 3).Value = 3
$ws.Cells.Item(142, 4).Value = 95
$ws.Cells.Item(142, 5).Value = 95
$ws.Cells.Item(142, 8).Value = 4

# Row 143
$ws.Cells.Item(143, 1).Value = 'Madagascar'
$ws.Cells.Item(143, 2).Value = 193
$ws.Cells.Item(143, 4).Value = 101
$ws.Cells.Item(143, 5).Value = 92
$ws.Cells.Item(143, 6).Value = 1
$ws.Cells.Item(143, 8).Value = 0

# Row 173
$ws.Cells.Item(173, 4).Value = 29
$ws.Cells.Item(173, 5).Value = 15

# Row 192
$ws.Cells.Item(192, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(192, 4).Value = 18
$ws.Cells.Item(192, 8).Value = 0

# Row 193
$ws.Cells.Item(193, 1).Value = 'Belice'
$ws.Cells.Item(193, 4).Value = 16
$ws.Cells.Item(193, 8).Value = 2

# Row 205
$ws.Cells.Item(205, 1).Value = 'Montserrat'
$ws.Cells.Item(205, 4).Value = 7
$ws.Cells.Item(205, 6).Value = 1
$ws.Cells.Item(205, 8).Value = 1

# Row 206
$ws.Cells.Item(206, 1).Value = 'Seychelles'
$ws.Cells.Item(206, 4).Value = 8
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 8).Value = 0
